$d = $word.ActiveDocument

$rng = $d.Content
if ($rng.Find.Execute("Individual Consultation Guide ", $true, $false, $false, $false, $false, $true, 1, $false)) {
    $rng.Text = "Guia di Konsulta Individual "
} else {
    Write-Output "NOT FOUND: 0"
}
$rng = $d.Content
if ($rng.Find.Execute("Learning about the family, child, and relationship between parent and child.", $true, $false, $false, $false, $false, $true, 1, $false)) {
    $rng.Text = "Siña tokante e famia, yu i relashon entre mayor i yu."
} else {
    Write-Output "NOT FOUND: 1"
}
$rng = $d.Content
if ($rng.Find.Execute("You should use individual consultations as an opportunity to learn as much as you can about the family circumstances and challenges. This is important as when a parent begins to change the way that they interact with the children in the household, it will affect everyone else and can sometimes generate conflict. ", $true, $false, $false, $false, $false, $true, 1, $false)) {
    $rng.Text = "Bo mester usa konsultanan individual komo un oportunidat pa siña mas tantu ku bo por tokante e sirkunstansianan i retonan di famia. Esaki ta importante ya ku ora un mayor kuminsá kambia e manera ku e ta interkambiá ku e yunan den e kas, e lo afektá tur otro hende i tin bia por generá konflikto. "
} else {
    Write-Output "NOT FOUND: 2"
}
$rng = $d.Content
if ($rng.Find.Execute("This conversation will provide you with a greater understanding of some of the challenges that families are experiencing at home. It will also equip you with information that might be useful in helping parents to establish weekly goals and solve issues that arise when parents practice new skills at home.", $true, $false, $false, $false, $false, $true, 1, $false)) {
    $rng.Text = "E kombersashon akí lo duna bo un komprondementu mas grandi di algun di e retonan ku famianan ta eksperensiando na kas. E lo ekipá bo tambe ku informashon ku por ta útil pa yuda mayornan establesé metanan semanal i solushoná asuntunan ku ta surgi ora mayornan ta praktiká abilidatnan nobo na kas."
} else {
    Write-Output "NOT FOUND: 3"
}
$rng = $d.Content
if ($rng.Find.Execute("Identifying an overall goal for the programme.", $true, $false, $false, $false, $false, $true, 1, $false)) {
    $rng.Text = "Identifiká un meta general pa e programa."
} else {
    Write-Output "NOT FOUND: 4"
}
$rng = $d.Content
if ($rng.Find.Execute("It is also important that you help parents to identify ONE positive, specific, and realistic goal for their target child during the programme.", $true, $false, $false, $false, $false, $true, 1, $false)) {
    $rng.Text = "Tambe ta importante pa bo yuda mayornan identifiká UN meta positivo, spesífiko i realístiko pa e yu ku nan ta enfoká riba dje durante e programa."
} else {
    Write-Output "NOT FOUND: 5"
}
$rng = $d.Content
if ($rng.Find.Execute("To help parents set a goal, it is important that you ask parents to describe their expectations about the programme. In doing this, help them identify a specific goal about how they would like to change their relationship with their child or their child’s behaviour. ", $true, $false, $false, $false, $false, $true, 1, $false)) {
    $rng.Text = "Pa yuda mayornan pone un meta, ta importante pa bo pidi mayornan pa deskribí nan ekspektativanan tokante e programa. Ora bo ta hasi esaki, yuda nan identifiká un meta spesífiko tokante kon nan lo ke kambia nan relashon ku nan yu òf e komportashon di nan yu. "
} else {
    Write-Output "NOT FOUND: 6"
}
$rng = $d.Content
if ($rng.Find.Execute("The goal should be specific", $true, $false, $false, $false, $false, $true, 1, $false)) {
    $rng.Text = "E meta mester ta spesífiko"
} else {
    Write-Output "NOT FOUND: 7"
}
$rng = $d.Content
if ($rng.Find.Execute("Parents will often give vague goals such as, “I want my child to be good,” or “I want my child to succeed at school.” You may need to guide them to be more specific about what the parent means by “good” or “succeed at school.” ", $true, $false, $false, $false, $false, $true, 1, $false)) {
    $rng.Text = "Hopi bia mayornan lo duna metanan vago manera, “Mi ke pa mi yu ta bon mucha,” òf “Mi ke pa mi yu tin éksito na skol.” Kisas bo mester guia nan pa ta mas spesífiko tokante kiko e mayor ke men ku “bon mucha” òf “tin éksito na skol.” "
} else {
    Write-Output "NOT FOUND: 8"
}
$rng = $d.Content
if ($rng.Find.Execute("It is your job to help parents describe a behaviour as specifically as possible. You will need to use questions to get them to be more and more specific. We recommend that you ask parents what “being good” or “succeeding at school” means to them. Further, exploring when, where, and why they would like this behaviour to happen will help make their goal more specific.", $true, $false, $false, $false, $false, $true, 1, $false)) {
    $rng.Text = "Ta bo trabou pa yuda mayornan deskribí un komportashon mas spesífiko posibel. Lo bo mester usa pregunta pa logra pa nan ta mas i mas spesífiko. Nos ta rekomendá pa bo puntra mayornan kiko “ta bon mucha” òf “tin éksito na skol” ta nifiká pa nan. Ademas, eksplorando ki ora, unda i dikon nan lo ke pa e komportashon aki sosodé, lo yuda hasi nan meta mas spesífiko."
} else {
    Write-Output "NOT FOUND: 9"
}
$rng = $d.Content
if ($rng.Find.Execute("Making goals specific will help the parents identify behaviours that they can actually help encourage to happen more often as they develop more positive relationships with their child. ", $true, $false, $false, $false, $false, $true, 1, $false)) {
    $rng.Text = "Hasiendo metanan spesífiko lo yuda e mayornan identifiká komportashonnan ku nan por yuda enkurashá pa sosodé mas frekuente segun ku nan ta desaroyá relashonnan mas positivo ku nan yu. "
} else {
    Write-Output "NOT FOUND: 10"
}
$rng = $d.Content
if ($rng.Find.Execute("A parent who wants her child “to be good” may eventually have a goal such as “I want my child to greet his elders in a respectful way when they enter the house.” ", $true, $false, $false, $false, $false, $true, 1, $false)) {
    $rng.Text = "Un mayor ku ke pa su yu “ta bon mucha” eventualmente por tin un meta manera “Mi ke pa mi yu kuminda su grandinan na un manera respetuoso ora e drenta kas.” "
} else {
    Write-Output "NOT FOUND: 11"
}
$rng = $d.Content
if ($rng.Find.Execute("Likewise, you can help a parent who wants his child to “succeed at school” to identify a more specific goal such as, “I want my child to do her homework immediately after coming home from school.” ", $true, $false, $false, $false, $false, $true, 1, $false)) {
    $rng.Text = "Di mes manera, bo por yuda un mayor ku ke pa su yu “tin éksito na skol” pa identifiká un meta mas spesífiko manera, “Mi ke pa mi yu hasi su hùiswèrk mesora despues di yega kas for di skol.” "
} else {
    Write-Output "NOT FOUND: 12"
}
$rng = $d.Content
if ($rng.Find.Execute("The goal should also be stated in a positive way.", $true, $false, $false, $false, $false, $true, 1, $false)) {
    $rng.Text = "Mester deklará e meta tambe na un manera positivo."
} else {
    Write-Output "NOT FOUND: 13"
}
$rng = $d.Content
if ($rng.Find.Execute("For example, instead of saying, “I want my child to stop swearing at me,” a parent should be helped to state the behaviour that s/he wants to see: “I want my child to use friendly words when talking to me.” ", $true, $false, $false, $false, $false, $true, 1, $false)) {
    $rng.Text = "Por ehèmpel, en bes di bisa, “Mi ke pa mi yu stòp di papia palabra mahos ku mi,” mester yuda un mayor pa deklará e komportashon ku e ke mira: “Mi ke pa mi yu usa palabranan amistoso ora e ta papia ku mi.” "
} else {
    Write-Output "NOT FOUND: 14"
}
$rng = $d.Content
if ($rng.Find.Execute("When you and the parent are happy with the specific, positive, and realistic goal, you should write it down and send it to them as a text message for safekeeping.", $true, $false, $false, $false, $false, $true, 1, $false)) {
    $rng.Text = "Ora abo i e mayor ta kontentu ku e meta spesífiko, positivo i realístiko, bo mester skibié i mand’é p'é komo un mensahe di teksto pa warda esaki."
} else {
    Write-Output "NOT FOUND: 15"
}
$rng = $d.Content
if ($rng.Find.Execute("The goal should also be realistic", $true, $false, $false, $false, $false, $true, 1, $false)) {
    $rng.Text = "E meta tambe mester ta realístiko"
} else {
    Write-Output "NOT FOUND: 16"
}
$rng = $d.Content
if ($rng.Find.Execute("Parents will often want to choose goals that are not achievable during the programme or sometimes are impossible for their child’s developmental stage. You can gently guide parents to have more realistic expectations for their children.", $true, $false, $false, $false, $false, $true, 1, $false)) {
    $rng.Text = "Mayornan hopi bia lo ke skohe metanan ku no ta alkansabel durante e programa òf tin bia ta imposibel pa e fase di desaroyo di nan yu. Bo por guia mayornan suavemente pa tin espektativaanan mas realístiko pa nan yunan."
} else {
    Write-Output "NOT FOUND: 17"
}
$rng = $d.Content
if ($rng.Find.Execute("For example, if a parent says that she wants her child to pass her school examinations at the end of the year, you might want to ask her what specific behaviours her child needs to do in order to make that possible.", $true, $false, $false, $false, $false, $true, 1, $false)) {
    $rng.Text = "Por ehèmpel, si un mayor bisa ku e ke pa su yu pasa su èksamennan di skol na fin di aña, bo por puntr’é ki komportashonnan spesífiko su yu mester hasi pa hasi esei posibel."
} else {
    Write-Output "NOT FOUND: 18"
}
$rng = $d.Content
if ($rng.Find.Execute("Likewise, if a parent wants his 2-year-old child to be able to get dressed by himself in the morning, you need to explore whether that is a reasonable expectation and why the child cannot get dressed by himself. You then could help the parent identify a more realistic goal that would help his child develop the skills to get dressed by himself such as, “I would like my child to cooperate with me while I help him get dressed in the morning.” ", $true, $false, $false, $false, $false, $true, 1, $false)) {
    $rng.Text = "Di mes manera, si un mayor ke pa su yu di 2 aña por bisti su so mainta, bo mester eksplorá si esei ta un ekspektativa rasonabel i dikon e yu no por bisti su so. Despues bo por yuda e mayor identifiká un meta mas realístiko ku lo yuda su yu desaroyá e abilidatnan pa bisti su so manera, “Mi lo ke pa mi yu koperá ku mi miéntras mi ta yud’é bisti mainta.” "
} else {
    Write-Output "NOT FOUND: 19"
}
$rng = $d.Content
if ($rng.Find.Execute("Discuss logistics and technology for the group sessions", $true, $false, $false, $false, $false, $true, 1, $false)) {
    $rng.Text = "Diskutí logístika i teknologia pa e seshonnan di grupo"
} else {
    Write-Output "NOT FOUND: 20"
}
$rng = $d.Content
if ($rng.Find.Execute("Finally, individual consultations are opportunities to discuss any logistical matters with the parents about the time for the group sessions, access to a cellphone, any data load/airtime needs, and technological questions.", $true, $false, $false, $false, $false, $true, 1, $false)) {
    $rng.Text = "Finalmente, konsultanan individual ta oportunidatnan pa papia tokante kualke asuntu logístiko ku e mayornan tokante e ora pa e seshonnan di grupo, akseso na un telefòn selular, kualke nesesidat di rekargá data i preguntanan teknológiko."
} else {
    Write-Output "NOT FOUND: 21"
}
$rng = $d.Content
if ($rng.Find.Execute("Suggested Structure for Individual Consultations before Group Sessions:", $true, $false, $false, $false, $false, $true, 1, $false)) {
    $rng.Text = "Struktura Sugerí pa Konsultanan Individual promé ku Seshonnan di Grupo:"
} else {
    Write-Output "NOT FOUND: 22"
}
$rng = $d.Content
if ($rng.Find.Execute("A. Introduce yourselves to the parent and the whole family if present.", $true, $false, $false, $false, $false, $true, 1, $false)) {
    $rng.Text = "A. Introdusí boso mes na e mayor i henter famia si ta presente."
} else {
    Write-Output "NOT FOUND: 23"
}
$rng = $d.Content
if ($rng.Find.Execute("B. Provide an overview of programme (This can just be the basics such as how the programme will help them with their own and their child’s behaviour)", $true, $false, $false, $false, $false, $true, 1, $false)) {
    $rng.Text = "B. Duna un bista general di e programa (Esaki por ta djis e kosnan básiko manera kon e programa lo yuda nan ku nan mes i nan yu su komportashon)"
} else {
    Write-Output "NOT FOUND: 24"
}
$rng = $d.Content
if ($rng.Find.Execute("C. Ask about the family environment at home: ", $true, $false, $false, $false, $false, $true, 1, $false)) {
    $rng.Text = "C. Puntra tokante e ambiente di famia na kas: "
} else {
    Write-Output "NOT FOUND: 25"
}
$rng = $d.Content
if ($rng.Find.Execute("What is happening at home?", $true, $false, $false, $false, $false, $true, 1, $false)) {
    $rng.Text = "Kiko ta pasando na kas?"
} else {
    Write-Output "NOT FOUND: 26"
}
$rng = $d.Content
if ($rng.Find.Execute("Who else lives there? How many children are there? Husband/Wife? Partner? Grandparents?", $true, $false, $false, $false, $false, $true, 1, $false)) {
    $rng.Text = "Ken mas ta biba einan? Kuantu mucha tin? Esposo/Esposa? Pareha? Grandinan?"
} else {
    Write-Output "NOT FOUND: 27"
}
$rng = $d.Content
if ($rng.Find.Execute("Who else provides care for children?", $true, $false, $false, $false, $false, $true, 1, $false)) {
    $rng.Text = "Ken mas ta duna kuido na mucha?"
} else {
    Write-Output "NOT FOUND: 28"
}
$rng = $d.Content
if ($rng.Find.Execute("What kind of support do you already receive from close friends and family members that you can trust nearby?", $true, $false, $false, $false, $false, $true, 1, $false)) {
    $rng.Text = "Ki tipo di sosten bo ta risibí kaba di amigunan i miembronan di famia yegá ku bo por konfia den serkania?"
} else {
    Write-Output "NOT FOUND: 29"
}
$rng = $d.Content
if ($rng.Find.Execute("D. Discuss with the parent about his/her relationship with their child:", $true, $false, $false, $false, $false, $true, 1, $false)) {
    $rng.Text = "D. Kombersá ku e mayor tokante su relashon ku su yu:"
} else {
    Write-Output "NOT FOUND: 30"
}
$rng = $d.Content
if ($rng.Find.Execute("If the programme is being delivered as part of a study: ", $true, $false, $false, $false, $false, $true, 1, $false)) {
    $rng.Text = "Si ta implementando e programa komo parti di un estudio: "
} else {
    Write-Output "NOT FOUND: 31"
}
$rng = $d.Content
if ($rng.Find.Execute("Remind the parent that they will be focusing on the target child selected during the assessment.", $true, $false, $false, $false, $false, $true, 1, $false)) {
    $rng.Text = "Kòrda e mayor ku e lo ta enfoká riba e yu ku el a skohe pa enfoká riba dje durante e evaluashon."
} else {
    Write-Output "NOT FOUND: 32"
}
$rng = $d.Content
if ($rng.Find.Execute("If the parent has not selected a specific child to focus on during the programme: ", $true, $false, $false, $false, $false, $true, 1, $false)) {
    $rng.Text = "Si e mayor no a selektá un yu spesífiko pa enfoká riba dje durante e programa: "
} else {
    Write-Output "NOT FOUND: 33"
}
$rng = $d.Content
if ($rng.Find.Execute("Ask the parent to select one child to focus on during the programme. This child should be between the ages of 2 and 17 years old. If the parent has more than one child between this age range, s/he should select the child with whom s/he is having the most difficult relationship or challenges when managing the child’s behaviour. ", $true, $false, $false, $false, $false, $true, 1, $false)) {
    $rng.Text = "Pidi e mayor pa selektá un yu pa enfoká riba dje durante e programa. E yu akí mester ta entre e edat di 2 i 17 aña. Si e mayor tin mas ku un yu entre e rango di edat akí, e mester selektá e yu ku ken aworaki e tin e relashon òf retonan mas difísil ora di manehá e komportashon di e yu. "
} else {
    Write-Output "NOT FOUND: 34"
}
$rng = $d.Content
if ($rng.Find.Execute("You can also reassure the parent that the skills learned in the programme may be applicable to all of the other children in his/her family, but that s/he should focus on this one child during the group discussions and home practice.", $true, $false, $false, $false, $false, $true, 1, $false)) {
    $rng.Text = "Bo por sigurá e mayor tambe ku e abilidatnan ku ta siña den e programa por ta aplikabel pa tur e otro muchanan den su famia, pero ku e mester enfoká riba e mucha akí durante e diskushonnan di grupo i práktika na kas."
} else {
    Write-Output "NOT FOUND: 35"
}
$rng = $d.Content
if ($rng.Find.Execute("What is life like at home with your child?", $true, $false, $false, $false, $false, $true, 1, $false)) {
    $rng.Text = "Kon bida ta na kas ku bo yu?"
} else {
    Write-Output "NOT FOUND: 36"
}
$rng = $d.Content
if ($rng.Find.Execute("What is your relationship with your child like?", $true, $false, $false, $false, $false, $true, 1, $false)) {
    $rng.Text = "Kon bo relashon ku bo yu ta?"
} else {
    Write-Output "NOT FOUND: 37"
}
$rng = $d.Content
if ($rng.Find.Execute("What are some challenges that you are facing in terms of managing your child’s behaviour?", $true, $false, $false, $false, $false, $true, 1, $false)) {
    $rng.Text = "Kua ta algun reto ku bo ta konfrontando pa loke ta trata manehá e komportashon di bo yu?"
} else {
    Write-Output "NOT FOUND: 38"
}
$rng = $d.Content
if ($rng.Find.Execute("Are there other challenges that make it difficult to be a parent?", $true, $false, $false, $false, $false, $true, 1, $false)) {
    $rng.Text = "Tin otro retonan ku ta hasié difísil pa ta un mayor?"
} else {
    Write-Output "NOT FOUND: 39"
}
$rng = $d.Content
if ($rng.Find.Execute("E. Parent Goals for the Programme", $true, $false, $false, $false, $false, $true, 1, $false)) {
    $rng.Text = "E. Metanan di Mayor pa e Programa"
} else {
    Write-Output "NOT FOUND: 40"
}
$rng = $d.Content
if ($rng.Find.Execute("What are your goals, expectations, or hopes for you and your child in general and from this programme? ", $true, $false, $false, $false, $false, $true, 1, $false)) {
    $rng.Text = "Kiko ta bo metanan, ekspektativanan òf speransanan pa abo i bo yu en general i pa e programa akí? "
} else {
    Write-Output "NOT FOUND: 41"
}
$rng = $d.Content
if ($rng.Find.Execute("Help the parent identify ONE the specific, positive, and realistic goal.", $true, $false, $false, $false, $false, $true, 1, $false)) {
    $rng.Text = "Yuda e mayor identifiká UN e meta spesífiko, positivo i realístiko."
} else {
    Write-Output "NOT FOUND: 42"
}
$rng = $d.Content
if ($rng.Find.Execute("Write the parent’s goal down at the bottom of your participant/parent profile.", $true, $false, $false, $false, $false, $true, 1, $false)) {
    $rng.Text = "Skibi e meta di e mayor na parti abou di bo profil di partisipante/mayor."
} else {
    Write-Output "NOT FOUND: 43"
}
$rng = $d.Content
if ($rng.Find.Execute("F. Discuss practicalities", $true, $false, $false, $false, $false, $true, 1, $false)) {
    $rng.Text = "F. Diskutí kosnan práktiko"
} else {
    Write-Output "NOT FOUND: 44"
}
$rng = $d.Content
if ($rng.Find.Execute("Timing of the session", $true, $false, $false, $false, $false, $true, 1, $false)) {
    $rng.Text = "Tempu di e seshon"
} else {
    Write-Output "NOT FOUND: 45"
}
$rng = $d.Content
if ($rng.Find.Execute("Platform that will be used for ParentChat (and how to download the app if they do not already have it)", $true, $false, $false, $false, $false, $true, 1, $false)) {
    $rng.Text = "Plataforma ku lo usa pa Mayor Konektá (i kon pa baha e aplikahson si nan no tin esaki kaba)"
} else {
    Write-Output "NOT FOUND: 46"
}
$rng = $d.Content
if ($rng.Find.Execute("Other technological questions that are related to participating in the programme. What is the participant’s digital literacy? Do participants know how to open messages? Respond?", $true, $false, $false, $false, $false, $true, 1, $false)) {
    $rng.Text = "Otro preguntanan teknológiko ku ta relashoná ku partisipá na e programa. Kiko ta e alfabetisashon digital di e partisipante? Partisipantenan sa kon pa habri mensahenan? Kontestá?"
} else {
    Write-Output "NOT FOUND: 47"
}
$rng = $d.Content
if ($rng.Find.Execute("Clarify literacy level of parents – you may need to send audio messages instead of text messages if parents have difficulty reading. You may need to explain how a participant listens to an audio message.", $true, $false, $false, $false, $false, $true, 1, $false)) {
    $rng.Text = "Aklará e nivel di alfabetisashon di mayornan – lo bo mester manda mensahenan di oudio en bes di mensahenan di teksto si mayornan tin difikultat pa lesa. Lo bo por mester splika kon un partisipante ta skucha un mensahe di oudio."
} else {
    Write-Output "NOT FOUND: 48"
}
$rng = $d.Content
if ($rng.Find.Execute("G. Any other questions?", $true, $false, $false, $false, $false, $true, 1, $false)) {
    $rng.Text = "G. Tin otro pregunta?"
} else {
    Write-Output "NOT FOUND: 49"
}
